$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 102, pushing the existing rows 102-154 down to 104-156.
$ws.Range("A102:A103").EntireRow.Insert()

# Populate the two newly inserted rows with the new weekly price entries for Jengibre.
# Row 102 (new)
$ws.Cells.Item(102, 1).Value = 8
$ws.Cells.Item(102, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(102, 3).Value = "Coquimbo"
$ws.Cells.Item(102, 4).Value = 45126
$ws.Cells.Item(102, 5).Value = 4
$ws.Cells.Item(102, 6).Value = 100114007
$ws.Cells.Item(102, 7).Value = "Jengibre"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 400
$ws.Cells.Item(102, 11).Value = 17000
$ws.Cells.Item(102, 12).Value = 18000
$ws.Cells.Item(102, 13).Value = 17500
$ws.Cells.Item(102, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(102, 15).Value = "Perú"
$ws.Cells.Item(102, 16).Value = 1346
$ws.Cells.Item(102, 17).Value = 13
$ws.Cells.Item(102, 18).Value = "Hortaliza"

# Row 103 (new)
$ws.Cells.Item(103, 1).Value = 8
$ws.Cells.Item(103, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(103, 3).Value = "Coquimbo"
$ws.Cells.Item(103, 4).Value = 45126
$ws.Cells.Item(103, 5).Value = 4
$ws.Cells.Item(103, 6).Value = 100114007
$ws.Cells.Item(103, 7).Value = "Jengibre"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 440
$ws.Cells.Item(103, 11).Value = 17000
$ws.Cells.Item(103, 12).Value = 18000
$ws.Cells.Item(103, 13).Value = 17500
$ws.Cells.Item(103, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(103, 15).Value = "Perú"
$ws.Cells.Item(103, 16).Value = 1346
$ws.Cells.Item(103, 17).Value = 13
$ws.Cells.Item(103, 18).Value = "Hortaliza"
